$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C5").Value = 88.17582417582418
$ws.Range("C6").Value = 90.36170212765957
$ws.Range("C7").Value = 94.09999999999999
$ws.Range("C8").Value = 91.40449438202248
$ws.Range("C9").Value = 103.1702127659574
$ws.Range("C10").Value = 85.07368421052631
$ws.Range("C11").Value = 85.69230769230769
